# Update legend labels to "large-" prefixed dataset names
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "large-KDD99"
$ws.Range("C1").Value = "large-CoverType"
$ws.Range("D1").Value = "large-KDD98"

# Widen columns B, C, D to fit the new longer legend text
$ws.Columns.Item(2).ColumnWidth = 13
$ws.Columns.Item(3).ColumnWidth = 19
$ws.Columns.Item(4).ColumnWidth = 12.6

# Move the active selection to E10 (as recorded in the saved view state)
$ws.Range("E10").Select()
